$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / update the workbook title to reflect the new "as-of" date.
$ws.Name = "Through 2021-12-01"

# --- September row (row 11): 2021 arrest count revised down by one -------
$ws.Range("U11").Value = 170
$ws.Range("V11").Value = 0.0395

# --- November row (row 13): no longer a partial month, revise 2021 values
$ws.Range("A13").Value = "November"
$ws.Range("T13").Value = 5
$ws.Range("U13").Value = 197
$ws.Range("V13").Value = 0.0248

# --- Insert a new "December (through 12-01)" row ahead of the Total row -
$ws.Rows.Item(14).Insert()

# Clear whatever formatting Excel auto-propagated into the new blank row so
# we can rebuild only the cells that the new row actually has data for.
$ws.Range("B14:V14").Clear()

# Copy the label-cell formatting (bold font + border) from the row above.
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").Value = "December (through 12-01)"

$ws.Range("F14").Value = 3
$ws.Range("I14").Value = 4
$ws.Range("L14").Value = 3
$ws.Range("R14").Value = 4
$ws.Range("U14").Value = 6

# --- Total row (was row 14, now shifted to row 15): refresh the sums -----
$ws.Range("B15").Value = 33
$ws.Range("C15").Value = 258
$ws.Range("D15").Value = 0.1134
$ws.Range("E15").Value = 60
$ws.Range("F15").Value = 506
$ws.Range("G15").Value = 0.106
$ws.Range("H15").Value = 63
$ws.Range("I15").Value = 762
$ws.Range("J15").Value = 0.0764
$ws.Range("K15").Value = 74
$ws.Range("L15").Value = 611
$ws.Range("M15").Value = 0.108
$ws.Range("N15").Value = 54
$ws.Range("O15").Value = 480
$ws.Range("P15").Value = 0.1011
$ws.Range("Q15").Value = 64
$ws.Range("R15").Value = 1204
$ws.Range("S15").Value = 0.0505
$ws.Range("T15").Value = 97
$ws.Range("U15").Value = 1553
$ws.Range("V15").Value = 0.0588
